$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying species-occurrence records were refreshed/re-ordered.
# Row 2 now holds the record that used to be on row 3 (with an updated
# sort key and coordinates); row 3 now holds what used to be on row 5;
# row 4 now holds what used to be on row 2; row 5 now holds what used to
# be on row 4.

# Row 2
$ws.Range("A2").Value = 112165478
$ws.Range("B2").Value = 77650
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("L2").ClearContents()
$ws.Range("Q2").Value = 333707
$ws.Range("R2").Value = 6627024

# Row 3
$ws.Range("A3").Value = 112165405
$ws.Range("B3").Value = 77650
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 333617
$ws.Range("R3").Value = 6627003

# Row 4
$ws.Range("A4").Value = 112164902
$ws.Range("B4").Value = 96735
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("Q4").Value = 333235
$ws.Range("R4").Value = 6626921
$ws.Range("AC4").ClearContents()

# Row 5
$ws.Range("A5").Value = 112165178
$ws.Range("B5").Value = 95369
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 2389
$ws.Range("F5").Value = "Purpurmylia"
$ws.Range("G5").Value = "Mylia taylorii"
$ws.Range("H5").Value = "(Hook.) Gray"
$ws.Range("Q5").Value = 333515
$ws.Range("R5").Value = 6626887
$ws.Range("AC5").Value = "På murken låga"
